$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells. Each target cell is forced to
# keep its original (General) style while storing the new value as text,
# since the source data are inline strings (e.g. "29.401.93", "  -0.27%  ")
# that must not be auto-coerced into numbers by Excel.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '29.401.93'
Set-TextValue "D3" '1.846.31'
Set-TextValue "E3" '  -0.27%  '
Set-TextValue "D4" '0.9993'
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '239.62'
Set-TextValue "E5" '  -0.68%  '
Set-TextValue "D6" '0.6319'
Set-TextValue "E6" '  -0.26%  '
Set-TextValue "E7" '  +0.02%  '
Set-TextValue "D8" '0.07567'
Set-TextValue "E8" '  -0.07%  '
Set-TextValue "D9" '0.2932'
Set-TextValue "E9" '  -0.96%  '
Set-TextValue "D10" '24.54'
Set-TextValue "E10" '  -0.31%  '
Set-TextValue "D11" '0.07716'
Set-TextValue "E11" '  -0.14%  '
Set-TextValue "D12" '1.870.07'
Set-TextValue "E12" '  -5.79%  '
Set-TextValue "D13" '5.004'
Set-TextValue "E13" '  +0.13%  '
Set-TextValue "D14" '0.6803'
Set-TextValue "E14" '  -0.66%  '
Set-TextValue "D15" '0.00001045'
Set-TextValue "E15" '  +5.36%  '
Set-TextValue "E16" '  +0.51%  '
Set-TextValue "D17" '2.123.05'
Set-TextValue "E17" '  -6.24%  '
Set-TextValue "D18" '6.171'
Set-TextValue "E18" '  -0.19%  '
Set-TextValue "D19" '29.448.10'
Set-TextValue "E19" '  -0.03%  '
Set-TextValue "D20" '228.71'
Set-TextValue "E20" '  -1.32%  '
Set-TextValue "E21" '  -0.21%  '
Set-TextValue "D22" '1.000'
Set-TextValue "E22" '  +0.04%  '
Set-TextValue "E23" '  -1.57%  '
Set-TextValue "E24" '  +0.07%  '
Set-TextValue "D25" '156.70'
Set-TextValue "E25" '  +0.46%  '
Set-TextValue "D26" '0.1394'
Set-TextValue "E26" '  +0.45%  '
Set-TextValue "D27" '8.336'
Set-TextValue "E27" '  -0.83%  '
Set-TextValue "D28" '17.62'
Set-TextValue "E28" '  -0.45%  '
Set-TextValue "D29" '1.466'
Set-TextValue "E29" '  -0.21%  '
Set-TextValue "E30" '  +3.47%  '
Set-TextValue "E31" '  -1.95%  '
Set-TextValue "D32" '4.102'
Set-TextValue "E32" '  -0.66%  '
Set-TextValue "D33" '4.024'
Set-TextValue "E33" '  +0.16%  '
Set-TextValue "E34" '  -0.27%  '
Set-TextValue "D35" '1.157'
Set-TextValue "E35" '  -0.15%  '
Set-TextValue "D36" '0.7094'
Set-TextValue "E36" '  -0.98%  '
Set-TextValue "E37" '  -0.06%  '
Set-TextValue "D38" '1.251.02'
Set-TextValue "E38" '  -0.19%  '
Set-TextValue "E39" '  +0.21%  '
Set-TextValue "D40" '2.772'
Set-TextValue "E40" '  -1.07%  '
Set-TextValue "D41" '6.385'
Set-TextValue "E41" '  +4.79%  '
Set-TextValue "D42" '0.9032'
Set-TextValue "E42" '  -0.20%  '
Set-TextValue "E43" '  +0.03%  '
Set-TextValue "D44" '101.77'
Set-TextValue "E44" '  +0.06%  '
Set-TextValue "D45" '65.90'
Set-TextValue "E45" '  -1.57%  '
Set-TextValue "D46" '0.00000000120'
Set-TextValue "E46" '  +1.42%  '
Set-TextValue "D47" '7.098'
Set-TextValue "E47" '  -0.60%  '
Set-TextValue "D48" '0.4000'
Set-TextValue "E48" '  -0.54%  '
Set-TextValue "B49" 'EnergySwap'
Set-TextValue "C49" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D49" '8.941'
Set-TextValue "E49" '  -2.86%  '
Set-TextValue "B50" 'RenderToken'
Set-TextValue "C50" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D50" '1.675'
Set-TextValue "E50" '  -0.48%  '
Set-TextValue "D51" '0.1123'
Set-TextValue "E51" '  -0.14%  '
